# Auto-generated script applying the Masamune_Profits value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$cols = @("H","I","J","K","L","M","N")
$vals = @(25000784,28571610,5000,85714830,15000,-85714658,-15344)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "61").Value = $vals[$i] }
$cols = @("H","I","K","M")
$vals = @(33650.25,819.35297,819.35297,678.64703)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "98").Value = $vals[$i] }
$cols = @("H","I","K","M")
$vals = @(33650.25,819.35297,2458.05891,-8.058910000000196)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "122").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(1682.6875,1516,1738.25,13644,15644.25,-11184,-20564.25)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "125").Value = $vals[$i] }
$cols = @("H","I","K","M")
$vals = @(1675665,4525741.5,13577224.5,-13574674.5)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "137").Value = $vals[$i] }

$ws = $wb.Worksheets.Item("ARM")
$cols = @("H","I","J","K","L","M","N")
$vals = @(1759.6066,1313.5814,2825.111,1313.5814,2825.111,-1101.5814,-3249.111)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "61").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(1754.102,1590.6666,2091.1875,1590.6666,2091.1875,-716.6666,-3839.1875)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "74").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(1754.102,1590.6666,2091.1875,7953.333000000001,10455.9375,-3585.333000000001,-19191.9375)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "77").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(1624.7333,1519.4783,1970.5714,4558.4349,5911.7142,-2108.4349,-10811.7142)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "122").Value = $vals[$i] }
$cols = @("H","J","L","N")
$vals = @(39429,39429,39429,-49229)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "123").Value = $vals[$i] }
$cols = @("H","J","L","N")
$vals = @(52226.668,52226.668,52226.668,-62366.668)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "134").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(1759.6066,1313.5814,2825.111,3940.7442,8475.332999999999,-1390.7442,-13575.333)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "136").Value = $vals[$i] }

$ws = $wb.Worksheets.Item("BSM")
$cols = @("H","I","J","K","L","M","N")
$vals = @(367,240,1002,240,1002,-67,-1348)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "22").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(2716.8853,1899.2759,3457.8438,5697.8277,10373.5314,-3162.8277,-15443.5314)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "134").Value = $vals[$i] }

$ws = $wb.Worksheets.Item("CRP")
$cols = @("H","I","J","K","L","M","N")
$vals = @(7941812,1698.5,18528630,1698.5,18528630,-1403.5,-18529220)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "31").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(7941812,1698.5,18528630,1698.5,18528630,-1496.5,-18529034)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "34").Value = $vals[$i] }
$cols = @("H","J","L","N")
$vals = @(40000,40000,40000,-40588)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "52").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(1940.1621,1282.25,3154.7693,1282.25,3154.7693,-1079.25,-3560.7693)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "58").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(93176.53999999999,120831.1,994.6667,362493.3,2984.0001,-360043.3,-7884.0001)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "122").Value = $vals[$i] }
$cols = @("H","I","K","M")
$vals = @(63844,1345.7273,4037.1819,-1507.1819)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "132").Value = $vals[$i] }
$cols = @("H","I","K","M")
$vals = @(409175.8,486933.47,1460800.41,-1458265.41)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "134").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(1940.1621,1282.25,3154.7693,3846.75,9464.3079,-1296.75,-14564.3079)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "136").Value = $vals[$i] }

$ws = $wb.Worksheets.Item("CUL")
$cols = @("H","I","J","K","L","M","N")
$vals = @(2342.5715,1350,3666,4050,10998,-3780,-11538)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "64").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(2342.5715,1350,3666,4050,10998,-3114,-12870)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "67").Value = $vals[$i] }
$cols = @("H","I","K","M")
$vals = @(1004410,1500615,4501845,-4497007)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "120").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(2179.14,14702.714,1159.779,44108.142,3479.337,-39068.142,-13559.337)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "131").Value = $vals[$i] }

$ws = $wb.Worksheets.Item("GSM")
$cols = @("H","I","K","M")
$vals = @(1677.7778,1683.3334,5050.0002,-2600.0002)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "122").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(76933144,142873490,2754.6667,428620470,8264.000100000001,-428618000,-13204.0001)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "126").Value = $vals[$i] }

$ws = $wb.Worksheets.Item("LTW")
$cols = @("H","I","J","K","L","M")
$vals = @(3231.9,3231.9,0,3231.9,0,-3095.9)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "40").Value = $vals[$i] }
$ws.Range("N40").ClearContents()
$cols = @("H","I","J","K","L","M","N")
$vals = @(68960.47,85550.664,2599.6667,256651.992,7799.000100000001,-254201.992,-12699.0001)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "122").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(3607.682,2498.1333,5985.2856,7494.3999,17955.8568,-4964.3999,-23015.8568)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "132").Value = $vals[$i] }
$cols = @("H","I","K","M")
$vals = @(2637.3333,1913.6923,5741.0769,-3191.0769)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "136").Value = $vals[$i] }

$ws = $wb.Worksheets.Item("WVR")
$cols = @("H","I","J","K","L","M")
$vals = @(1588931.5,1588931.5,0,4766794.5,0,-4764344.5)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "122").Value = $vals[$i] }
$ws.Range("N122").ClearContents()
$cols = @("H","J","L","N")
$vals = @(30000,30000,30000,-39800)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "123").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(1892579.4,3954021.5,2924.1667,11862064.5,8772.500100000001,-11859534.5,-13832.5001)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "132").Value = $vals[$i] }
$cols = @("H","I","J","K","L","M","N")
$vals = @(615053.1,865049.3,1426.091,2595147.9,4278.272999999999,-2592597.9,-9378.272999999999)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "136").Value = $vals[$i] }
